# Weekly update: insert two new price records (week of 2022-06-02 / Excel
# serial 44714) above the existing history for "Vega Monumental Concepción -
# Limón", pushing all the older rows (formerly 458-490) down to 460-492.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 458:459; everything currently at row 458 and
# below shifts down by two rows (row 490 -> 492), matching the new
# dimension A1:T492.
$ws.Rows("458:459").Insert()

# New row 458: Limón "1a amarillo", Provincia de Melipilla
$ws.Cells.Item(458, 1).Value = 11
$ws.Cells.Item(458, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(458, 3).Value = "Bíobío"
$ws.Cells.Item(458, 4).Value = 44714
$ws.Cells.Item(458, 5).Value = 8
$ws.Cells.Item(458, 6).Value = "Fruta"
$ws.Cells.Item(458, 7).Value = 100102
$ws.Cells.Item(458, 8).Value = "Cítricos"
$ws.Cells.Item(458, 9).Value = 100102003
$ws.Cells.Item(458, 10).Value = "Limón"
$ws.Cells.Item(458, 11).Value = "Sin especificar"
$ws.Cells.Item(458, 12).Value = "1a amarillo"
$ws.Cells.Item(458, 13).Value = 300
$ws.Cells.Item(458, 14).Value = 10000
$ws.Cells.Item(458, 15).Value = 10000
$ws.Cells.Item(458, 16).Value = 10000
$ws.Cells.Item(458, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(458, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(458, 19).Value = 625
$ws.Cells.Item(458, 20).Value = 16

# New row 459: Limón "2a amarillo", Provincia de Melipilla
$ws.Cells.Item(459, 1).Value = 11
$ws.Cells.Item(459, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(459, 3).Value = "Bíobío"
$ws.Cells.Item(459, 4).Value = 44714
$ws.Cells.Item(459, 5).Value = 8
$ws.Cells.Item(459, 6).Value = "Fruta"
$ws.Cells.Item(459, 7).Value = 100102
$ws.Cells.Item(459, 8).Value = "Cítricos"
$ws.Cells.Item(459, 9).Value = 100102003
$ws.Cells.Item(459, 10).Value = "Limón"
$ws.Cells.Item(459, 11).Value = "Sin especificar"
$ws.Cells.Item(459, 12).Value = "2a amarillo"
$ws.Cells.Item(459, 13).Value = 300
$ws.Cells.Item(459, 14).Value = 8000
$ws.Cells.Item(459, 15).Value = 8000
$ws.Cells.Item(459, 16).Value = 8000
$ws.Cells.Item(459, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(459, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(459, 19).Value = 500
$ws.Cells.Item(459, 20).Value = 16
